# Add working set of sequences
# Replace the 192 German verb entries in column B (word) of rows 2..193
# with the new working-set ordering, leaving every other column untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$words = @(
    'äußern', 'feuern', 'seufzen', 'stoßen', 'fügen', 'achten', 'sprengen', 'malen',
    'fühlen', 'kümmern', 'irren', 'sitzen', 'wachsen', 'runden', 'schämen', 'streichen',
    'ärgern', 'wetten', 'eignen', 'schwören', 'werden', 'warnen', 'hauen', 'süßen',
    'feiern', 'geben', 'stürzen', 'gründen', 'sagen', 'drehen', 'lehnen', 'wenden',
    'rauben', 'stellen', 'weichen', 'lesen', 'jagen', 'schultern', 'bellen', 'morden',
    'erben', 'hoffen', 'wirken', 'flehen', 'jubeln', 'stammen', 'leugnen', 'rasen',
    'pfeifen', 'stehlen', 'arten', 'bremsen', 'dienen', 'schrecken', 'freuen', 'decken',
    'schreiten', 'rechnen', 'suchen', 'hören', 'schenken', 'starren', 'graben', 'folgen',
    'weinen', 'räumen', 'knarren', 'kosten', 'garen', 'mauern', 'orten', 'greifen',
    'ahnen', 'werfen', 'gelten', 'bilden', 'ehren', 'schleppen', 'spinnen', 'sterben',
    'ändern', 'grüßen', 'quälen', 'drohen', 'trösten', 'ruhen', 'schlucken', 'zünden',
    'loben', 'plaudern', 'binden', 'fließen', 'platzen', 'wellen', 'lockern', 'pissen',
    'betteln', 'klettern', 'holen', 'lügen', 'stürmen', 'wehen', 'pressen', 'trauen',
    'kämpfen', 'fehlen', 'liefern', 'wehtun', 'kratzen', 'laufen', 'schnellen', 'töten',
    'schwingen', 'fischen', 'brauchen', 'bluten', 'fahren', 'tauchen', 'schulden', 'rufen',
    'scheiden', 'scheitern', 'boxen', 'altern', 'bergen', 'fallen', 'wehren', 'retten',
    'treiben', 'münzen', 'bitten', 'helfen', 'schwächen', 'proben', 'heilen', 'formen',
    'scheinen', 'machen', 'treten', 'zeigen', 'saufen', 'fällen', 'mögen', 'flüchten',
    'kaufen', 'spielen', 'klingen', 'siegen', 'schmecken', 'rühren', 'schalten', 'sparen',
    'spüren', 'testen', 'stecken', 'pflanzen', 'zögern', 'hacken', 'wundern', 'sorgen',
    'kehren', 'heulen', 'pflegen', 'tollen', 'filmen', 'wüten', 'streifen', 'enden',
    'zielen', 'tropfen', 'führen', 'albern', 'lassen', 'biegen', 'bieten', 'fangen',
    'kichern', 'sichern', 'liegen', 'trennen', 'husten', 'sinken', 'dringen', 'reizen',
    'messen', 'leisten', 'sperren', 'backen', 'planen', 'bauen', 'heben', 'merken'
)

for ($i = 0; $i -lt $words.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $words[$i]
}
